$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the header / data values (and drop the old "finalUrl" header text) ---
$ws.Range("A1").Value = "Path"
$ws.Range("A2").Value = "/basic_page/applied-biosystems-model-7000-real-time-thermalcycler/"
$ws.Range("A3").Value = "/basic_page/105-description/"
$ws.Range("A4").Value = "/basic_page/119-description/"
$ws.Range("A5").Value = "/basic_page/209-description/"
$ws.Range("A6").Value = "/basic_page/216-description/"

# --- 2. Remove the hyperlinks that used to live on A2:A4 ---
$ws.Range("A2").Hyperlinks.Delete()
$ws.Range("A3").Hyperlinks.Delete()
$ws.Range("A4").Hyperlinks.Delete()

# --- 3. A2:A4 go back to the plain default look (they used to carry the Hyperlink style) ---
$ws.Range("A2:A4").Style = "Normal"

# --- 4. Give A1 the new bold/white-on-blue header look ---
$header = $ws.Range("A1")
$header.Font.Bold = $true
$header.Font.Size = 14
$header.Font.Name = "Arial"
$header.Font.ThemeColor = 2
$header.Interior.PatternColor = 15983311
$header.Interior.ThemeColor = 5
$bottomBorder = $header.Borders.Item(9)
$bottomBorder.Color = 0
$bottomBorder.LineStyle = 1
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4108
$header.WrapText = $true

# --- 5. Narrow column A a bit ---
$ws.Columns.Item(1).ColumnWidth = 79.16

# --- 6. New (empty) row 14 styled like a hyperlink cell ---
$ws.Range("A14").Style = "Hyperlink"

# --- 7. Restore the selection that was active when the file was last saved ---
$ws.Range("B26").Select()
